# Add AERA 2021 conference presentation and AERA Open in-press article,
# and refresh the 2021 course offering details (D/S specialization).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3: EDLD 652 Data visualization for EDS -> 2021 offering ---
$ws.Range("D3").Value = "2019/2020/2021 "
$ws.Range("E3").Value = "27553/27120/27056"
$ws.Range("H3").Value = "https://dataviz-2021.netlify.app/"
$ws.Range("I3").Value = "https://github.com/uo-datasci-specialization/c2-dataviz-2021"

# --- Row 4: EDLD 653 Functional programming for EDS -> 2021 offering ---
$ws.Range("D4").Value = "2019/2020/2021"
$ws.Range("E4").Value = "35699/32066/36713"
$ws.Range("H4").Value = "https://fp-2021.netlify.app/"
$ws.Range("I4").Value = "https://github.com/uo-datasci-specialization/c3-fp-2021"

# --- Row 6: Capstone, now actually scheduled for 2021 (was "planned") ---
$ws.Range("D6").Value = 2021
$ws.Range("E6").Value = 27140

# --- Row 7: was a placeholder "Hierarchical Linear Modeling II [planned]" row,
#     now the AERA '21 conference presentation + AERA Open in-press article ---
$ws.Range("C7").Value = "Spring"
$ws.Range("D7").Value = 2021
$ws.Range("E7").Value = 36724
$ws.Range("H7").Value = "https://mlm2-2021.netlify.app/"
$ws.Range("I7").Value = "https://github.com/datalorax/mlm2"

# Drop the now-stale hyperlink relationship (H3 keeps its Hyperlink cell style).
$ws.Hyperlinks.Delete()

# Leave the cursor where the author last left it.
$ws.Range("C7").Select()
